$wb = $excel.ActiveWorkbook

# Sheet1: "Submit orders" - add rows 77-78
$ws1 = $wb.Worksheets.Item("Submit orders")
$ws1.Range("A77").Value = "10.06.2022 12:19 (Kyiv+Israel) 09:19 (UTC) 18:19 (Japan) 14:49 (India)"
$ws1.Range("B77").Value = 1.346
$ws1.Range("C77").Value = -0.5330000000000001
$ws1.Range("D77").Value = "***"
$ws1.Range("E77").Value = "***"

$ws1.Range("A78").Value = "10.07.2022 12:55 (Kyiv+Israel) 09:55 (UTC) 18:55 (Japan) 15:25 (India)"
$ws1.Range("B78").Value = 0.727
$ws1.Range("C78").Value = 0.08599999999999997
$ws1.Range("D78").Value = "***"
$ws1.Range("E78").Value = "***"

# Sheet2: "Submit internet survey" - add rows 72-73
$ws2 = $wb.Worksheets.Item("Submit internet survey")
$ws2.Range("A72").Value = "10.06.2022 15:18 (Kyiv+Israel) 12:18 (UTC) 21:18 (Japan) 17:48 (India)"
$ws2.Range("B72").Value = 177.221
$ws2.Range("C72").Value = -176.59
$ws2.Range("D72").Value = "***"
$ws2.Range("E72").Value = "***"

$ws2.Range("A73").Value = "10.07.2022 12:58 (Kyiv+Israel) 09:58 (UTC) 18:58 (Japan) 15:28 (India)"
$ws2.Range("B73").Value = 0.848
$ws2.Range("C73").Value = -0.217
$ws2.Range("D73").Value = "***"
$ws2.Range("E73").Value = "***"

# Sheet3: "Submit a phone survey" - add row 65
$ws3 = $wb.Worksheets.Item("Submit a phone survey")
$ws3.Range("A65").Value = "10.07.2022 13:47 (Kyiv+Israel) 10:47 (UTC) 19:47 (Japan) 16:17 (India)"
$ws3.Range("B65").Value = 1.676
$ws3.Range("C65").Value = -0.212
$ws3.Range("D65").Value = "***"
$ws3.Range("E65").Value = "***"

# Sheet4: "Checkertificate" - add row 75
$ws4 = $wb.Worksheets.Item("Checkertificate")
$ws4.Range("A75").Value = "10.07.2022 13:05 (Kyiv+Israel) 10:05 (UTC) 19:05 (Japan) 15:35 (India)"
$ws4.Range("B75").Value = 0.665
$ws4.Range("C75").Value = 0.02899999999999991
$ws4.Range("D75").Value = "***"
$ws4.Range("E75").Value = "***"
